$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume/1h (E) columns.
# D-column values are entered with a leading apostrophe so numeric-looking
# strings (e.g. "212.08") stay text instead of becoming Doubles, then the
# cell style is reset to Normal so no quote-prefix formatting lingers.

$ws.Range("D2").Value = "'26.296.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "'1.619.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("D5").Value = "'212.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "'0.0614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "'18.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.72%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'1.846.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "'1.623.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "'0.517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'26.305.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "'62.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("D18").Value = "'0.0₃0727"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'201.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "'4.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "'9.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("D25").Value = "'143.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'0.120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "'6.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("E30").Value = "  +10.71%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").Value = "'1.179.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.11%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "'0.806"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "'0.788"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").Value = "'5.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.74%  "
$ws.Range("D44").Value = "'1.757.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("D45").Value = "'93.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +15.03%  "
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").Value = "'53.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").Value = "  -0.27%  "
